# Apply the "added all unique SPZ to pocty prujezdu" edit.
#
# Summary of the change:
#  - Rename the "Pocet SPZ" column header (R1) to
#    "Celkovy pocet unikatnich SPZ " (note trailing space).
#  - The R column used to hold a per-row formula (P/A). The author
#    replaced that with a single, manually-entered "total unique SPZ"
#    number in R3, highlighted like the other manually-entered rows
#    (red text on a pale-yellow fill), and removed the now-unused R
#    formulas from every other row (R4:R16, R18).
#  - The conditional-format rule that highlights B21:O21 when the
#    difference is exactly 0 changes its fill color.
#  - The active selection moves to R17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column header text -------------------------------------------------
$ws.Range("R1").Value = "Celkový počet unikátních SPZ "

# --- 2. Replace the per-row formula in R with a single literal total ------
# R3 becomes a hard-coded "unique SPZ" count instead of the old =P3/A3
# formula, and gets flagged with the same red-on-yellow style used
# elsewhere in the sheet for manually supplied numbers.
$r3 = $ws.Range("R3")
$r3.Value = 32761
$r3.Font.Bold = $false
$r3.Font.Color = 255        # RGB(FF,00,00) red
$r3.Interior.Color = 13431551   # RGB(FF,F2,CC) pale yellow
$r3.HorizontalAlignment = -4108   # xlCenter
$r3.VerticalAlignment = -4108     # xlCenter
$r3.WrapText = $true

# Remove the now-obsolete R formulas on every other data / total row.
$ws.Range("R4:R16").Clear()
$ws.Range("R18").Clear()

# --- 3. Conditional formatting tweak on B21:O21 ----------------------------
# The "equal to 0" rule's highlight fill changes color.
$diffRange = $ws.Range("B21:O21")
for ($i = 1; $i -le $diffRange.FormatConditions.Count; $i++) {
    $fc = $diffRange.FormatConditions.Item($i)
    if ($fc.Operator -eq 3) {
        # xlEqual
        $fc.Interior.Color = 12379352   # RGB(D8,E4,BC)
    }
}

# --- 4. Row heights tweak up slightly (header + the wrapped target row) ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(20).RowHeight = 30

# --- 5. Leave the selection where the author left it -----------------------
$ws.Range("R17").Select()

Write-Host "edit applied"
